$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Category" in A1, matching the style of the other header cells (B1:W1)
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) # xlPasteFormats

# Remove the header-style formatting from A2:A46 (data cells, not headers)
$ws.Range("A2:A46").ClearFormats()
